$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ligand/receptor TPM-derived metrics (NATMI lrc2p output) for Epha4-Efnb1 pairs.
# Only the cells whose values change with the new TPM input are touched; row/column
# layout (sender/target cluster labels in A:D) is unchanged.
$updates = @{
    "2" = @{ "E"=3; "F"=1; "G"=8.081040666666667; "H"=24.243122; "I"=0.4661250698616886; "J"=0.4661250698616886; "M"=13.16594766666667; "N"=39.497843; "O"=0.6940777873489595; "P"=0.6940777873489595; "Q"=106.3945585095385; "R"=957.5510265858461; "S"=0.32352705711748; "T"=0.32352705711748 }
    "3" = @{ "E"=3; "F"=1; "G"=8.081040666666667; "H"=24.243122; "I"=0.4661250698616886; "J"=0.4661250698616886; "O"=0.1706596770095176; "P"=0.1706596770095176; "Q"=26.16026808776; "R"=235.44241278984; "S"=0.07954875386863461; "T"=0.07954875386863461 }
    "4" = @{ "E"=3; "F"=1; "G"=8.081040666666667; "H"=24.243122; "I"=0.4661250698616886; "J"=0.4661250698616886; "N"=7.697376999999999; "O"=0.1352625356415228; "P"=0.1352625356415228; "Q"=20.73427218788822; "R"=186.608449690994; "S"=0.06304925887557394; "T"=0.06304925887557396 }
    "5" = @{ "I"=0.4037865631294714; "J"=0.4037865631294715; "M"=13.16594766666667; "N"=39.497843; "O"=0.6940777873489595; "P"=0.6940777873489595; "Q"=92.16559222826514; "R"=829.4903300543862; "S"=0.2802592842981445; "T"=0.2802592842981445 }
    "6" = @{ "I"=0.4037865631294714; "J"=0.4037865631294715; "O"=0.1706596770095176; "P"=0.1706596770095176; "S"=0.06891008444445879; "T"=0.0689100844444588 }
    "7" = @{ "I"=0.4037865631294714; "J"=0.4037865631294715; "N"=7.697376999999999; "O"=0.1352625356415228; "P"=0.1352625356415228; "S"=0.05461719438686811; "T"=0.05461719438686813 }
    "8" = @{ "I"=0.1300883670088399; "J"=0.1300883670088399; "M"=13.16594766666667; "N"=39.497843; "O"=0.6940777873489595; "P"=0.6940777873489595; "Q"=29.69309155424578; "R"=267.237823988212; "S"=0.09029144593333499; "T"=0.09029144593333499 }
    "9" = @{ "I"=0.1300883670088399; "J"=0.1300883670088399; "O"=0.1706596770095176; "P"=0.1706596770095176; "S"=0.02220083869642421; "T"=0.02220083869642421 }
    "10" = @{ "I"=0.1300883670088399; "J"=0.1300883670088399; "N"=7.697376999999999; "O"=0.1352625356415228; "P"=0.1352625356415228; "Q"=5.786617765140887; "R"=52.07955988626799; "S"=0.0175960823790807; "T"=0.01759608237908071 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}